$d = $word.ActiveDocument

# 1. Expand on the underwater vehicle vision example and start a new sentence
#    before "Image restoration".
$r1 = $d.Content
$r1.Find.Execute(
    "underwater monitoring, underwater robots, etc, image restoration",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "underwater monitoring with the use of underwater vehicle vision, etc. Image restoration",
    2)

# 2. British spelling: color -> colour.
$r2 = $d.Content
$r2.Find.Execute(
    "image color contrast",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "image colour contrast",
    2)

# 3. Replace "Hence, we provide a residual learning-based" with
#    " We provide a machine learning-based" (note the extra space before
#    "We" left after the removed "Hence,").
$r3 = $d.Content
$r3.Find.Execute(
    "framework. Hence, we provide a residual learning-based",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "framework.  We provide a machine learning-based",
    2)
